# Auto-generated Excel COM-interop edit script
# Applies odds/data updates described by the commit diff
# to Jogos_do_Dia_Betfair_Back_Lay_2025-10-08.xlsx (Sheet1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("N2").Value = 2.84
$ws.Range("Q2").Value = 2.3
$ws.Range("AK2").Value = 170

# --- Row 3 ---
$ws.Range("F3").Value = 8.199999999999999
$ws.Range("G3").Value = 13
$ws.Range("I3").Value = 1.54
$ws.Range("J3").Value = 4.3
$ws.Range("K3").Value = 4.7
$ws.Range("L3").Value = 1.43
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 3.1
$ws.Range("P3").Value = 1.73
$ws.Range("X3").Value = 15
$ws.Range("AA3").Value = 14
$ws.Range("AF3").Value = 95
$ws.Range("AK3").Value = 220
$ws.Range("AL3").Value = 200
$ws.Range("AO3").Value = 11

# --- Row 4 ---
$ws.Range("G4").Value = 4.2
$ws.Range("I4").Value = 2.52
$ws.Range("J4").Value = 3.15
$ws.Range("L4").Value = 1.46
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 3.05
$ws.Range("O4").Value = 1.39
$ws.Range("P4").Value = 1.72
$ws.Range("Q4").Value = 2.12
$ws.Range("R4").Value = 1.27
$ws.Range("S4").Value = 3.9
$ws.Range("T4").Value = 1.84
$ws.Range("U4").Value = 1.96
$ws.Range("V4").Value = 1.66
$ws.Range("W4").Value = 1.33
$ws.Range("X4").Value = 12.5
$ws.Range("Y4").Value = 9.4
$ws.Range("Z4").Value = 15.5
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 12.5
$ws.Range("AC4").Value = 8
$ws.Range("AD4").Value = 12
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 26
$ws.Range("AG4").Value = 16
$ws.Range("AH4").Value = 20
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 1000
$ws.Range("AK4").Value = 50
$ws.Range("AL4").Value = 1000
$ws.Range("AM4").Value = 150
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000

# --- Row 5 ---
$ws.Range("F5").Value = 1.63
$ws.Range("G5").Value = 1.65
$ws.Range("K5").Value = 4.1
$ws.Range("L5").Value = 1.47
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 3.15
$ws.Range("O5").Value = 1.39
$ws.Range("P5").Value = 1.73
$ws.Range("R5").Value = 1.26
$ws.Range("S5").Value = 4
$ws.Range("T5").Value = 1.96
$ws.Range("U5").Value = 1.62
$ws.Range("V5").Value = 1.15
$ws.Range("W5").Value = 2.46
$ws.Range("X5").Value = 12.5
$ws.Range("Y5").Value = 19.5
$ws.Range("Z5").Value = 60
$ws.Range("AA5").Value = 250
$ws.Range("AB5").Value = 8.199999999999999
$ws.Range("AC5").Value = 9.4
$ws.Range("AD5").Value = 28
$ws.Range("AE5").Value = 140
$ws.Range("AF5").Value = 9.4
$ws.Range("AG5").Value = 11
$ws.Range("AH5").Value = 28
$ws.Range("AI5").Value = 160
$ws.Range("AJ5").Value = 17.5
$ws.Range("AK5").Value = 21
$ws.Range("AL5").Value = 50
$ws.Range("AM5").Value = 240
$ws.Range("AN5").Value = 13
$ws.Range("AO5").Value = 220

# --- Row 6 ---
$ws.Range("D6").Value = "Cuiaba"
$ws.Range("E6").Value = "Novorizontino"
$ws.Range("F6").Value = 2.12
$ws.Range("G6").Value = 2.74
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 4.7
$ws.Range("J6").Value = 2.36
$ws.Range("K6").Value = 4.1
$ws.Range("L6").Value = 1.61
$ws.Range("M6").Value = 1.01
$ws.Range("N6").Value = 1.46
$ws.Range("O6").Value = 1.01
$ws.Range("P6").Value = 1.46
$ws.Range("Q6").Value = 2.58
$ws.Range("R6").Value = 1.17
$ws.Range("S6").Value = 5.5
$ws.Range("T6").Value = 1.01
$ws.Range("U6").Value = 1.01
$ws.Range("V6").Value = 1.27
$ws.Range("W6").Value = 1.57
$ws.Range("X6").Value = 10
$ws.Range("Y6").Value = 11.5
$ws.Range("Z6").Value = 38
$ws.Range("AA6").Value = 1000
$ws.Range("AB6").Value = 9.199999999999999
$ws.Range("AC6").Value = 9.800000000000001
$ws.Range("AD6").Value = 25
$ws.Range("AE6").Value = 1000
$ws.Range("AF6").Value = 18.5
$ws.Range("AG6").Value = 18
$ws.Range("AH6").Value = 36
$ws.Range("AI6").Value = 1000
$ws.Range("AJ6").Value = 50
$ws.Range("AK6").Value = 50
$ws.Range("AL6").Value = 1000
$ws.Range("AM6").Value = 1000
$ws.Range("AN6").Value = 1000
$ws.Range("AO6").Value = 1000

# --- Row 7 ---
$ws.Range("D7").Value = "Operario PR"
$ws.Range("E7").Value = "Athletic Club"
$ws.Range("F7").Value = 1.71
$ws.Range("G7").Value = 2.16
$ws.Range("H7").Value = 4.6
$ws.Range("I7").Value = 6.8
$ws.Range("J7").Value = 2.56
$ws.Range("K7").Value = 3.7
$ws.Range("L7").Value = 1.56
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 1.25
$ws.Range("O7").Value = 1.51
$ws.Range("Q7").Value = 2.56
$ws.Range("R7").Value = 1.17
$ws.Range("S7").Value = 4.7
$ws.Range("T7").Value = 1.01
$ws.Range("U7").Value = 1.01
$ws.Range("V7").Value = 1.17
$ws.Range("W7").Value = 1.86
$ws.Range("X7").Value = 10
$ws.Range("Y7").Value = 19
$ws.Range("Z7").Value = 55
$ws.Range("AA7").Value = 1000
$ws.Range("AB7").Value = 8.800000000000001
$ws.Range("AC7").Value = 10.5
$ws.Range("AD7").Value = 32
$ws.Range("AE7").Value = 1000
$ws.Range("AF7").Value = 13.5
$ws.Range("AG7").Value = 16
$ws.Range("AH7").Value = 34
$ws.Range("AI7").Value = 1000
$ws.Range("AJ7").Value = 30
$ws.Range("AK7").Value = 32
$ws.Range("AL7").Value = 75
$ws.Range("AM7").Value = 1000
$ws.Range("AN7").Value = 1000
$ws.Range("AO7").Value = 1000
